$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# --- Header row relabelling (row 1) ---
# "Study id" -> "id"
$ws.Range("C1").Value = "id"
# "RR (t2 vs. t1)" -> "effect"
$ws.Range("F1").Value = "effect"
# "SE" -> "se"
$ws.Range("I1").Value = "se"

# --- New effect-size (se) formulas, column I, rows 2-9 ---
# I2 and I3 entered individually, I4:I9 entered as one range (creates a
# shared formula group spanning I4:I9, matching how it was filled down)
$ws.Range("I2").Formula = "=(H2-G2)/(1.96^2)"
$ws.Range("I3").Formula = "=(H3-G3)/(1.96^2)"
$ws.Range("I4:I9").Formula = "=(H4-G4)/(1.96^2)"

# --- New label cell K8, styled like F1 (which now reads "effect") ---
$ws.Range("F1").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("K8").Value = "RR (t2 vs. t1)"

# --- Update current selection to reflect last-edited cell ---
[void]$ws.Range("I12").Select()
